$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.028.12"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.70"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.14"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5023"
$ws.Range("E7").Value = "  -2.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4032"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08231"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.413"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.908.73"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.286"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.20"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06492"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.068.10"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.192"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.30"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.130.41"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.25"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.293"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.97"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.131"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.011"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.818"
$ws.Range("E34").Value = "  +4.89%  "
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.349"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06433"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2163"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.905"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.205"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6442"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.38"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.217"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.41"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.197"
$ws.Range("E46").Value = "  +6.30%  "
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.635"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.47"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.214"
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.79"
$ws.Range("E51").Value = "  -0.73%  "
